# 977: Add GS to extract process and GS tab to example files
#
# Adds a new "GS" worksheet (as the last tab, after "CMS") containing a
# single header row with the columns used by the GS extract:
#   Contact_ID, Contact_Date, Contact_Type_Code, Contact_Type_Desc,
#   OM_Name, OM_Key, OM_Grade, OM_Team_Key, OM_Provider_Code

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

$gsSheet = $wb.Worksheets.Add($null, $lastSheet)
$gsSheet.Name = "GS"

$gsSheet.Range("A1").Value = "Contact_ID"
$gsSheet.Range("B1").Value = "Contact_Date"
$gsSheet.Range("C1").Value = "Contact_Type_Code"
$gsSheet.Range("D1").Value = "Contact_Type_Desc"
$gsSheet.Range("E1").Value = "OM_Name"
$gsSheet.Range("F1").Value = "OM_Key"
$gsSheet.Range("G1").Value = "OM_Grade"
$gsSheet.Range("H1").Value = "OM_Team_Key"
$gsSheet.Range("I1").Value = "OM_Provider_Code"

$headerRange = $gsSheet.Range("A1:I1")
$headerRange.Font.Color = 0

$gsSheet.Activate()
$headerRange.Select() | Out-Null
